# Insert a new price-record row above row 50 (Macroferia Regional de Talca - Piña,
# Caramelo / Segunda), shifting the existing rows 50-151 down to 51-152.
# The new row captures a price observation dated serial 44469 (2021-09-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(50).Insert()

$ws.Cells.Item(50, 1).Value  = 5
$ws.Cells.Item(50, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(50, 3).Value  = "Maule"
$ws.Cells.Item(50, 4).Value  = 44469
$ws.Cells.Item(50, 5).Value  = 7
$ws.Cells.Item(50, 6).Value  = "Fruta"
$ws.Cells.Item(50, 7).Value  = 100108
$ws.Cells.Item(50, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(50, 9).Value  = 100108005
$ws.Cells.Item(50, 10).Value = "Piña"
$ws.Cells.Item(50, 11).Value = "Caramelo"
$ws.Cells.Item(50, 12).Value = "Segunda"
$ws.Cells.Item(50, 13).Value = 100
$ws.Cells.Item(50, 14).Value = 20000
$ws.Cells.Item(50, 15).Value = 20000
$ws.Cells.Item(50, 16).Value = 20000
$ws.Cells.Item(50, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(50, 18).Value = "Ecuador"
$ws.Cells.Item(50, 19).Value = 1429
$ws.Cells.Item(50, 20).Value = 14
